# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Each affected cell currently holds the Excel serial date 45189
# (2023-09-20) and must be bumped forward by one day to 45190
# (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
